$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# The leave card grew two new entries inside the "Undertime" block of the
# table (between the existing "SL(1-0-0)"/"FL(5-0-0)"/"VL(21-0-0)" rows and
# the "2023" year marker). This is modeled as two single-row inserts:
#   - a new row at (old) row 38, holding a UT(0-4-0) / 0.5-day entry
#   - a new row at (old) row 40 -> becomes row 41 after the first insert,
#     holding a new UT(0-2-1) / 0.252-day entry
# xlShiftDown = -4121
$ws.Rows("38:38").Insert(-4121)
$ws.Rows("41:41").Insert(-4121)

# The table range needs to explicitly grow to cover the two extra rows
# (header row 8 .. new last row 134).
$tbl.Resize($ws.Range("A8:K134"))

# Seed the two brand-new rows from existing "typical" data rows so they pick
# up the exact same cell styles used throughout the table (rather than the
# engine minting fresh style records for them).
$ws.Range("A37:K37").Copy($ws.Range("A38:K38"))
$ws.Range("A36:K36").Copy($ws.Range("A41:K41"))

$earnedFormula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Row 38: UT(0-4-0), Absence Undertime W/O Pay = 0.5
$ws.Range("A38").ClearContents()
$ws.Range("B38").Value = "UT(0-4-0)"
$ws.Range("C38").ClearContents()
$ws.Range("D38").Value = 0.5
$ws.Range("G38").Formula = $earnedFormula
$ws.Range("H38").ClearContents()
$ws.Range("K38").ClearContents()

# Row 41: UT(0-2-1), Absence Undertime W/O Pay = 0.252
$ws.Range("A41").ClearContents()
$ws.Range("B41").Value = "UT(0-2-1)"
$ws.Range("C41").ClearContents()
$ws.Range("D41").Value = 0.252
$ws.Range("G41").Formula = $earnedFormula
$ws.Range("H41").ClearContents()
$ws.Range("K41").ClearContents()

# The row-insert can leave the trailing (pre-existing) table rows with a
# stale "[@EARNED]" structured reference that now points outside the table;
# reassert the calculated-column formula on them so they match every other
# row again.
$ws.Range("G133").Formula = $earnedFormula
$ws.Range("G134").Formula = $earnedFormula

# CONVERTION sheet: the "late calculator" lookup key moved from 2 to 4,
# which ripples into G3 via SUMIFS.
$wsConv = $wb.Worksheets.Item("CONVERTION")
$wsConv.Range("E3").Value = 4

# Restore the recorded cursor position in the frozen (bottomLeft) pane.
$null = $ws.Range("F39").Select()
